# Daily attendance processing - reorder "Recorded By" (column G) entries.
# For every data row, the comma-separated list of recorders in column G is
# reversed in order (e.g. "a, b, c" -> "c, b, a").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $value = $cell.Value2

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = $value -split ",\s*"

    if ($parts.Count -gt 1) {
        $revParts = $parts[($parts.Count - 1)..0]
        $cell.Value2 = [string]::Join(", ", $revParts)
    }
}
